$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @{ E=3; G=282.7879796666667; H=848.363939; I=0.9674521741401267; J=0.9674521741401266; K=3; M=117.044563; N=351.133689; O=0.3245365645427815; P=0.3245365645427815; Q=33098.79550173788; R=297889.1595156409; S=0.3139736049548815; T=0.3139736049548815 }
    3  = @{ E=3; G=282.7879796666667; H=848.363939; I=0.9674521741401267; J=0.9674521741401266; K=3; M=101.5800373333333; N=304.740112; O=0.281657135515876;  P=0.281657135515876;  Q=28725.61353195791; R=258530.5217876212; S=0.2724898081169145; T=0.2724898081169145 }
    4  = @{ E=3; G=282.7879796666667; H=848.363939; I=0.9674521741401267; J=0.9674521741401266; K=3; M=142.0267893333333; N=426.080368; O=0.3938062999413425; P=0.3938062999413425; Q=40163.46881411662; R=361471.2193270496; S=0.3809887610683307; T=0.3809887610683306 }
    5  = @{ E=3; G=7.714696666666668;  H=23.14409;   I=0.02639291836872237; J=0.02639291836872237; K=3; M=117.044563; N=351.133689; O=0.3245365645427815; P=0.3245365645427815; Q=902.9633000275568; R=8126.66970024801; S=0.008565467055643232; T=0.00856546705564323 }
    6  = @{ E=3; G=7.714696666666668;  H=23.14409;   I=0.02639291836872237; J=0.02639291836872237; K=3; M=101.5800373333333; N=304.740112; O=0.281657135515876;  P=0.281657135515876;  Q=783.6591754153424; R=7052.932578738081; S=0.007433753785638691; T=0.007433753785638688 }
    7  = @{ E=3; G=7.714696666666668;  H=23.14409;   I=0.02639291836872237; J=0.02639291836872237; K=3; M=142.0267893333333; N=426.080368; O=0.3938062999413425; P=0.3938062999413425; Q=1095.693598247236; R=9861.242384225121; S=0.01039369752744045; T=0.01039369752744045 }
    8  = @{ E=3; G=1.799090333333333;  H=5.397271;   I=0.006154907491150983; J=0.006154907491150983; K=3; M=117.044563; N=351.133689; O=0.3245365645427815; P=0.3245365645427815; Q=210.5737418625243; R=1895.163676762719; S=0.00199749253225677; T=0.00199749253225677 }
    9  = @{ E=3; G=1.799090333333333;  H=5.397271;   I=0.006154907491150983; J=0.006154907491150983; K=3; M=101.5800373333333; N=304.740112; O=0.281657135515876;  P=0.281657135515876;  Q=182.7516632260391; R=1644.764969034352; S=0.001733573613322793; T=0.001733573613322792 }
    10 = @{ E=3; G=1.799090333333333;  H=5.397271;   I=0.006154907491150983; J=0.006154907491150983; K=3; M=142.0267893333333; N=426.080368; O=0.3938062999413425; P=0.3938062999413425; Q=255.5190237639698; R=2299.671213875728; S=0.00242384134557142; T=0.00242384134557142 }
}

foreach ($rowNum in $data.Keys) {
    $row = $data[$rowNum]
    foreach ($col in $row.Keys) {
        $ws.Range("$col$rowNum").Value = $row[$col]
    }
}
